# Apply edit: rebuild rows 2-7 with cumulative conversation content,
# and remove rows 8-22 (conversation_17486_processed.xlsx)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = '[{"role": "assistant", "content": "Chào cậu! Bữa nay thầy giáo robot trên sao HỎA dạy tớ cách luyện nói theo nhịp như… hát rap đó! Mỗi cụm từ giống như một đoạn nhạc nhỏ, có nhịp điệu riêng."}, {"role": "assistant", "content": "Cứ theo nhịp là nói được tiếng Anh hay hơn liền! Giờ tụi mình thử cùng luyện nha. Cậu đã sẵn sàng chưa?"}, {"role": "user", "content": "Are you ready?"}]'
$ws.Range("C2").Value = 'Chào cậu! Bữa nay thầy giáo robot trên sao HỎA dạy tớ cách luyện nói theo nhịp như… hát rap đó! Mỗi cụm từ giống như một đoạn nhạc nhỏ, có nhịp điệu riêng.'

# Row 3
$ws.Range("A3").Value = '[{"role": "assistant", "content": "Chào cậu! Bữa nay thầy giáo robot trên sao HỎA dạy tớ cách luyện nói theo nhịp như… hát rap đó! Mỗi cụm từ giống như một đoạn nhạc nhỏ, có nhịp điệu riêng."}, {"role": "assistant", "content": "Cứ theo nhịp là nói được tiếng Anh hay hơn liền! Giờ tụi mình thử cùng luyện nha. Cậu đã sẵn sàng chưa?"}, {"role": "user", "content": "Are you ready?"}, {"role": "assistant", "content": "Great! Let''s start! chúng mình sẽ học nói tiếng Anh theo nhịp giống như hát nha! Bài hát này có nhiều cụm từ rất quen mà cũng rất vui. Bài đầu tiên nha. Đầu tiên, cậu hãy nghe thử cả bài nhé, đừng hát theo vội. Mình chỉ cần cảm nhịp thôi."}, {"role": "assistant", "content": "Giờ mình cùng lặp lại đoạn đầu tiên nha, câu chào cực kỳ thân thiện luôn nè."}, {"role": "assistant", "content": "Cậu nói lại cùng tớ nha."}, {"role": "user", "content": "Hello, hello, how are you? Hello, hello, hello, how are you?"}]'
$ws.Range("C3").Value = 'Cứ theo nhịp là nói được tiếng Anh hay hơn liền! Giờ tụi mình thử cùng luyện nha. Cậu đã sẵn sàng chưa?'

# Row 4
$ws.Range("A4").Value = '[{"role": "assistant", "content": "Chào cậu! Bữa nay thầy giáo robot trên sao HỎA dạy tớ cách luyện nói theo nhịp như… hát rap đó! Mỗi cụm từ giống như một đoạn nhạc nhỏ, có nhịp điệu riêng."}, {"role": "assistant", "content": "Cứ theo nhịp là nói được tiếng Anh hay hơn liền! Giờ tụi mình thử cùng luyện nha. Cậu đã sẵn sàng chưa?"}, {"role": "user", "content": "Are you ready?"}, {"role": "assistant", "content": "Great! Let''s start! chúng mình sẽ học nói tiếng Anh theo nhịp giống như hát nha! Bài hát này có nhiều cụm từ rất quen mà cũng rất vui. Bài đầu tiên nha. Đầu tiên, cậu hãy nghe thử cả bài nhé, đừng hát theo vội. Mình chỉ cần cảm nhịp thôi."}, {"role": "assistant", "content": "Giờ mình cùng lặp lại đoạn đầu tiên nha, câu chào cực kỳ thân thiện luôn nè."}, {"role": "assistant", "content": "Cậu nói lại cùng tớ nha."}, {"role": "user", "content": "Hello, hello, how are you? Hello, hello, hello, how are you?"}, {"role": "assistant", "content": "Nice! You did it! Tiếp theo là khi mình cảm thấy thật tuyệt, cùng nói theo nhịp nha!"}, {"role": "assistant", "content": "Nào, nói lại với tớ!"}, {"role": "user", "content": "I''m good. I''m great. I''m wonderful. I''m good. I''m great, great, great. I''m wonderful."}]'
$ws.Range("C4").Value = "Great! Let's start! chúng mình sẽ học nói tiếng Anh theo nhịp giống như hát nha! Bài hát này có nhiều cụm từ rất quen mà cũng rất vui. Bài đầu tiên nha. Đầu tiên, cậu hãy nghe thử cả bài nhé, đừng hát theo vội. Mình chỉ cần cảm nhịp thôi. "

# Row 5
$ws.Range("A5").Value = '[{"role": "assistant", "content": "Chào cậu! Bữa nay thầy giáo robot trên sao HỎA dạy tớ cách luyện nói theo nhịp như… hát rap đó! Mỗi cụm từ giống như một đoạn nhạc nhỏ, có nhịp điệu riêng."}, {"role": "assistant", "content": "Cứ theo nhịp là nói được tiếng Anh hay hơn liền! Giờ tụi mình thử cùng luyện nha. Cậu đã sẵn sàng chưa?"}, {"role": "user", "content": "Are you ready?"}, {"role": "assistant", "content": "Great! Let''s start! chúng mình sẽ học nói tiếng Anh theo nhịp giống như hát nha! Bài hát này có nhiều cụm từ rất quen mà cũng rất vui. Bài đầu tiên nha. Đầu tiên, cậu hãy nghe thử cả bài nhé, đừng hát theo vội. Mình chỉ cần cảm nhịp thôi."}, {"role": "assistant", "content": "Giờ mình cùng lặp lại đoạn đầu tiên nha, câu chào cực kỳ thân thiện luôn nè."}, {"role": "assistant", "content": "Cậu nói lại cùng tớ nha."}, {"role": "user", "content": "Hello, hello, how are you? Hello, hello, hello, how are you?"}, {"role": "assistant", "content": "Nice! You did it! Tiếp theo là khi mình cảm thấy thật tuyệt, cùng nói theo nhịp nha!"}, {"role": "assistant", "content": "Nào, nói lại với tớ!"}, {"role": "user", "content": "I''m good. I''m great. I''m wonderful. I''m good. I''m great, great, great. I''m wonderful."}, {"role": "assistant", "content": "Awesome! Cậu nói theo nhịp rất tốt! Câu chào quen thuộc quay lại nè! Mình cùng nói lại một lần nữa nhé!"}, {"role": "assistant", "content": "Cậu nhớ nhịp rồi đúng không? Cùng nói nào!"}, {"role": "user", "content": "hello hello how are you hello hello hello how are you"}]'
$ws.Range("C5").Value = 'Giờ mình cùng lặp lại đoạn đầu tiên nha, câu chào cực kỳ thân thiện luôn nè.'

# Row 6
$ws.Range("A6").Value = '[{"role": "assistant", "content": "Chào cậu! Bữa nay thầy giáo robot trên sao HỎA dạy tớ cách luyện nói theo nhịp như… hát rap đó! Mỗi cụm từ giống như một đoạn nhạc nhỏ, có nhịp điệu riêng."}, {"role": "assistant", "content": "Cứ theo nhịp là nói được tiếng Anh hay hơn liền! Giờ tụi mình thử cùng luyện nha. Cậu đã sẵn sàng chưa?"}, {"role": "user", "content": "Are you ready?"}, {"role": "assistant", "content": "Great! Let''s start! chúng mình sẽ học nói tiếng Anh theo nhịp giống như hát nha! Bài hát này có nhiều cụm từ rất quen mà cũng rất vui. Bài đầu tiên nha. Đầu tiên, cậu hãy nghe thử cả bài nhé, đừng hát theo vội. Mình chỉ cần cảm nhịp thôi."}, {"role": "assistant", "content": "Giờ mình cùng lặp lại đoạn đầu tiên nha, câu chào cực kỳ thân thiện luôn nè."}, {"role": "assistant", "content": "Cậu nói lại cùng tớ nha."}, {"role": "user", "content": "Hello, hello, how are you? Hello, hello, hello, how are you?"}, {"role": "assistant", "content": "Nice! You did it! Tiếp theo là khi mình cảm thấy thật tuyệt, cùng nói theo nhịp nha!"}, {"role": "assistant", "content": "Nào, nói lại với tớ!"}, {"role": "user", "content": "I''m good. I''m great. I''m wonderful. I''m good. I''m great, great, great. I''m wonderful."}, {"role": "assistant", "content": "Awesome! Cậu nói theo nhịp rất tốt! Câu chào quen thuộc quay lại nè! Mình cùng nói lại một lần nữa nhé!"}, {"role": "assistant", "content": "Cậu nhớ nhịp rồi đúng không? Cùng nói nào!"}, {"role": "user", "content": "hello hello how are you hello hello hello how are you"}, {"role": "assistant", "content": "Good job! Cậu nhớ rất giỏi! Đôi khi mình cũng cảm thấy hơi mệt hoặc đói. Vậy thì mình sẽ nói thế này nè"}, {"role": "assistant", "content": "Giờ thì cậu nhại theo tớ nha!"}, {"role": "user", "content": "I''m tired I''m hungry I''m not so good I''m tired I''m hungry I''m not so good"}]'
$ws.Range("C6").Value = 'Cậu nói lại cùng tớ nha.'

# Row 7
$ws.Range("A7").Value = '[{"role": "assistant", "content": "Chào cậu! Bữa nay thầy giáo robot trên sao HỎA dạy tớ cách luyện nói theo nhịp như… hát rap đó! Mỗi cụm từ giống như một đoạn nhạc nhỏ, có nhịp điệu riêng."}, {"role": "assistant", "content": "Cứ theo nhịp là nói được tiếng Anh hay hơn liền! Giờ tụi mình thử cùng luyện nha. Cậu đã sẵn sàng chưa?"}, {"role": "user", "content": "Are you ready?"}, {"role": "assistant", "content": "Great! Let''s start! chúng mình sẽ học nói tiếng Anh theo nhịp giống như hát nha! Bài hát này có nhiều cụm từ rất quen mà cũng rất vui. Bài đầu tiên nha. Đầu tiên, cậu hãy nghe thử cả bài nhé, đừng hát theo vội. Mình chỉ cần cảm nhịp thôi."}, {"role": "assistant", "content": "Giờ mình cùng lặp lại đoạn đầu tiên nha, câu chào cực kỳ thân thiện luôn nè."}, {"role": "assistant", "content": "Cậu nói lại cùng tớ nha."}, {"role": "user", "content": "Hello, hello, how are you? Hello, hello, hello, how are you?"}, {"role": "assistant", "content": "Nice! You did it! Tiếp theo là khi mình cảm thấy thật tuyệt, cùng nói theo nhịp nha!"}, {"role": "assistant", "content": "Nào, nói lại với tớ!"}, {"role": "user", "content": "I''m good. I''m great. I''m wonderful. I''m good. I''m great, great, great. I''m wonderful."}, {"role": "assistant", "content": "Awesome! Cậu nói theo nhịp rất tốt! Câu chào quen thuộc quay lại nè! Mình cùng nói lại một lần nữa nhé!"}, {"role": "assistant", "content": "Cậu nhớ nhịp rồi đúng không? Cùng nói nào!"}, {"role": "user", "content": "hello hello how are you hello hello hello how are you"}, {"role": "assistant", "content": "Good job! Cậu nhớ rất giỏi! Đôi khi mình cũng cảm thấy hơi mệt hoặc đói. Vậy thì mình sẽ nói thế này nè"}, {"role": "assistant", "content": "Giờ thì cậu nhại theo tớ nha!"}, {"role": "user", "content": "I''m tired I''m hungry I''m not so good I''m tired I''m hungry I''m not so good"}, {"role": "assistant", "content": "Great! Bây giờ mình cùng nói câu chào cuối nhé! Và cuối cùng, bài hát kết thúc bằng câu chào quen thuộc một lần nữa. Nói lại lần cuối nha!"}, {"role": "assistant", "content": "Cùng nói lại với tớ nào!"}, {"role": "user", "content": "hello hello hello how are you hello hello hello how are you"}]'
$ws.Range("C7").Value = "Nice! You did it! Tiếp theo là khi mình cảm thấy thật tuyệt, cùng nói theo nhịp nha! "

# Remove rows 8-22 entirely, so the used range/dimension becomes A1:C7
$ws.Range("A8:C22").Delete(-4162)
